$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @(376, "2021-08-22", "overview", "K02000001", "United Kingdom", 6492906, 32253,  49, 131640),
  @(377, "2021-08-23", "overview", "K02000001", "United Kingdom", 6524581, 31914,  40, 131680),
  @(378, "2021-08-24", "overview", "K02000001", "United Kingdom", 6555200, 30838, 174, 131854),
  @(379, "2021-08-25", "overview", "K02000001", "United Kingdom", 6590747, 35847, 149, 132003),
  @(380, "2021-08-26", "overview", "K02000001", "United Kingdom", 6628709, 38281, 140, 132143),
  @(381, "2021-08-27", "overview", "K02000001", "United Kingdom", 6666399, 38046, 100, 132243),
  @(382, "2021-08-28", "overview", "K02000001", "United Kingdom", 6698486, 32406, 133, 132376)
)

foreach ($r in $rows) {
    $rowNum = $r[0]

    # Column A holds a date-like string that Excel would otherwise auto-convert
    # to a date serial number, so force text formatting while writing it, then
    # drop back to the default "Normal" style so no extra style is attached.
    $cellA = $ws.Cells.Item($rowNum, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $r[1]
    $cellA.Style = "Normal"

    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
    $ws.Cells.Item($rowNum, 5).Value = $r[5]
    $ws.Cells.Item($rowNum, 6).Value = $r[6]
    $ws.Cells.Item($rowNum, 7).Value = $r[7]
    $ws.Cells.Item($rowNum, 8).Value = $r[8]
}
